$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 625476.9399999999
$ws.Range("I28").Value = 769561.4
$ws.Range("K28").Value = 769561.4
$ws.Range("M28").Value = -769076.4

$ws.Range("H70").Value = 3530.4194
$ws.Range("I70").Value = 2209.3
$ws.Range("J70").Value = 4159.524
$ws.Range("K70").Value = 6627.900000000001
$ws.Range("L70").Value = 12478.572
$ws.Range("M70").Value = -6357.900000000001
$ws.Range("N70").Value = -13018.572

$ws.Range("H73").Value = 3530.4194
$ws.Range("I73").Value = 2209.3
$ws.Range("J73").Value = 4159.524
$ws.Range("K73").Value = 6627.900000000001
$ws.Range("L73").Value = 12478.572
$ws.Range("M73").Value = -5691.900000000001
$ws.Range("N73").Value = -14350.572

$ws.Range("H135").Value = 2294.7778
$ws.Range("I135").Value = 2294.7778
$ws.Range("K135").Value = 20653.0002
$ws.Range("M135").Value = -18118.0002

$ws.Range("H138").Value = 4434.6
$ws.Range("J138").Value = 5081.6665
$ws.Range("L138").Value = 15244.9995
$ws.Range("N138").Value = -25524.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2767.8125
$ws.Range("I2").Value = 1684.3846
$ws.Range("J2").Value = 7462.6665
$ws.Range("K2").Value = 1684.3846
$ws.Range("L2").Value = 7462.6665
$ws.Range("M2").Value = -1571.3846
$ws.Range("N2").Value = -7688.6665

$ws.Range("H4").Value = 428
$ws.Range("I4").Value = 460
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 460
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -344
$ws.Range("N4").Value = -532

$ws.Range("H94").Value = 79999
$ws.Range("J94").Value = 79999
$ws.Range("L94").Value = 79999
$ws.Range("N94").Value = -81801

$ws.Range("H102").Value = 8802.5
$ws.Range("I102").Value = 8802.5
$ws.Range("K102").Value = 8802.5
$ws.Range("M102").Value = -7180.5

$ws.Range("H112").Value = 21933
$ws.Range("J112").Value = 21933
$ws.Range("L112").Value = 21933
$ws.Range("N112").Value = -24887

$ws.Range("H116").Value = 2767.8125
$ws.Range("I116").Value = 1684.3846
$ws.Range("J116").Value = 7462.6665
$ws.Range("K116").Value = 1684.3846
$ws.Range("L116").Value = 7462.6665
$ws.Range("M116").Value = 609.6153999999999
$ws.Range("N116").Value = -12050.6665

$ws.Range("H122").Value = 4905.6787
$ws.Range("I122").Value = 5184.657
$ws.Range("K122").Value = 15553.971
$ws.Range("M122").Value = -13103.971

$ws.Range("H132").Value = 6055.482
$ws.Range("I132").Value = 6098.6104
$ws.Range("J132").Value = 5502
$ws.Range("K132").Value = 18295.8312
$ws.Range("L132").Value = 16506
$ws.Range("M132").Value = -15765.8312
$ws.Range("N132").Value = -21566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2767.8125
$ws.Range("I3").Value = 1684.3846
$ws.Range("J3").Value = 7462.6665
$ws.Range("K3").Value = 1684.3846
$ws.Range("L3").Value = 7462.6665
$ws.Range("M3").Value = -1570.3846
$ws.Range("N3").Value = -7690.6665

$ws.Range("H20").Value = 1787.9412
$ws.Range("I20").Value = 1631.6154
$ws.Range("K20").Value = 1631.6154
$ws.Range("M20").Value = -1384.6154

$ws.Range("H86").Value = 2542.8572
$ws.Range("I86").Value = 2235
$ws.Range("K86").Value = 2235
$ws.Range("M86").Value = -1112

$ws.Range("H89").Value = 2542.8572
$ws.Range("I89").Value = 2235
$ws.Range("K89").Value = 11175
$ws.Range("M89").Value = -5559

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2245.2104
$ws.Range("I16").Value = 3312.111
$ws.Range("J16").Value = 1285
$ws.Range("K16").Value = 3312.111
$ws.Range("L16").Value = 1285
$ws.Range("M16").Value = -3025.111
$ws.Range("N16").Value = -1859

$ws.Range("H22").Value = 5316.125
$ws.Range("I22").Value = 5288.3335
$ws.Range("J22").Value = 5399.5
$ws.Range("K22").Value = 5288.3335
$ws.Range("L22").Value = 5399.5
$ws.Range("M22").Value = -4938.3335
$ws.Range("N22").Value = -6099.5

$ws.Range("H31").Value = 3413.2856
$ws.Range("I31").Value = 2177.3125
$ws.Range("J31").Value = 4454.1055
$ws.Range("K31").Value = 2177.3125
$ws.Range("L31").Value = 4454.1055
$ws.Range("M31").Value = -1882.3125
$ws.Range("N31").Value = -5044.1055

$ws.Range("H34").Value = 3413.2856
$ws.Range("I34").Value = 2177.3125
$ws.Range("J34").Value = 4454.1055
$ws.Range("K34").Value = 2177.3125
$ws.Range("L34").Value = 4454.1055
$ws.Range("M34").Value = -1975.3125
$ws.Range("N34").Value = -4858.1055

$ws.Range("H58").Value = 4794.5
$ws.Range("I58").Value = 2657.8333
$ws.Range("K58").Value = 2657.8333
$ws.Range("M58").Value = -2454.8333

$ws.Range("H105").Value = 1651.1052
$ws.Range("I105").Value = 1750
$ws.Range("J105").Value = 1515.125
$ws.Range("K105").Value = 1750
$ws.Range("L105").Value = 1515.125
$ws.Range("M105").Value = -3
$ws.Range("N105").Value = -5009.125

$ws.Range("H113").Value = 2245.2104
$ws.Range("I113").Value = 3312.111
$ws.Range("J113").Value = 1285
$ws.Range("K113").Value = 3312.111
$ws.Range("L113").Value = 1285
$ws.Range("M113").Value = -1142.111
$ws.Range("N113").Value = -5625

$ws.Range("H132").Value = 2407.818
$ws.Range("I132").Value = 2070.6
$ws.Range("K132").Value = 6211.799999999999
$ws.Range("M132").Value = -3681.799999999999

$ws.Range("H136").Value = 4794.5
$ws.Range("I136").Value = 2657.8333
$ws.Range("K136").Value = 7973.499899999999
$ws.Range("M136").Value = -5423.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 121.6
$ws.Range("I38").Value = 138
$ws.Range("J38").Value = 83.333336
$ws.Range("K38").Value = 414
$ws.Range("L38").Value = 250.000008
$ws.Range("M38").Value = -67
$ws.Range("N38").Value = -944.000008

$ws.Range("H97").Value = 448.91306
$ws.Range("I97").Value = 305.75
$ws.Range("K97").Value = 917.25
$ws.Range("M97").Value = -421.25

$ws.Range("H113").Value = 4797.8
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 3497.25
$ws.Range("K113").Value = 30000
$ws.Range("L113").Value = 10491.75
$ws.Range("M113").Value = -27830
$ws.Range("N113").Value = -14831.75

$ws.Range("H128").Value = 694746.25
$ws.Range("I128").Value = 694746.25
$ws.Range("K128").Value = 2084238.75
$ws.Range("M128").Value = -2079258.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 3337.5386
$ws.Range("I122").Value = 2763.5557
$ws.Range("K122").Value = 8290.667099999999
$ws.Range("M122").Value = -5840.667099999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1651.3334
$ws.Range("I55").Value = 421.6
$ws.Range("J55").Value = 7800
$ws.Range("K55").Value = 421.6
$ws.Range("L55").Value = 7800
$ws.Range("M55").Value = -248.6
$ws.Range("N55").Value = -8146

$ws.Range("H61").Value = 2909.2
$ws.Range("I61").Value = 3091.2307
$ws.Range("K61").Value = 3091.2307
$ws.Range("M61").Value = -2889.2307

$ws.Range("H113").Value = 2909.2
$ws.Range("I113").Value = 3091.2307
$ws.Range("K113").Value = 3091.2307
$ws.Range("M113").Value = -921.2307000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 164283.42
$ws.Range("I62").Value = 164283.42
$ws.Range("K62").Value = 164283.42
$ws.Range("M62").Value = -163659.42

$ws.Range("H65").Value = 164283.42
$ws.Range("I65").Value = 164283.42
$ws.Range("K65").Value = 821417.1000000001
$ws.Range("M65").Value = -818297.1000000001

$ws.Range("H136").Value = 5888.25
$ws.Range("I136").Value = 4280.8
$ws.Range("J136").Value = 30000
$ws.Range("K136").Value = 30000
$ws.Range("L136").Value = 90000
$ws.Range("M136").Value = -10292.4
$ws.Range("N136").Value = -95100
